$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1, styled like the existing header row (copy style from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I and J, rows 2-14
$data = @(
    @(6, 7),
    @(4, 5),
    @(9, 9),
    @(6, 6),
    @(9, 9),
    @(8, 9),
    @(4, 6),
    @(8, 8),
    @(9, 9),
    @(10, 10),
    @(8, 8),
    @(6, 6),
    @(6, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
